# Regenerate the handback-status report timestamps.
# These cells hold the "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" (and the Overview's
# "Latest HO Xliff Generate Date") values as text, stored as shared
# strings. Update each cell that shows an older timestamp to the newer
# run's timestamp.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 11:10:30"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first file row.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 11:10:26"
$wsZhCn.Range("K2").Value = "2016-09-04 11:10:44"

# de-de sheet: Correspond Handoff Datetime (shared with Overview's value
# above) / Correspond Handback DateTime for the first file row.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 11:10:30"
$wsDeDe.Range("K2").Value = "2016-09-04 11:10:51"
